$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- N2: 8.5 -> 9 ---
$ws.Range("N2").Value = 9

# --- Row 3: several odds updated in place ---
$ws.Range("G3").Value  = 6.5
$ws.Range("H3").Value  = 4.33
$ws.Range("I3").Value  = 1.36
$ws.Range("J3").Value  = 6.5
$ws.Range("N3").Value  = 15
$ws.Range("U3").Value  = 1.8
$ws.Range("V3").Value  = 1.91
$ws.Range("Y3").Value  = 21
$ws.Range("AB3").Value = 41
$ws.Range("AD3").Value = 9
$ws.Range("AE3").Value = 17
$ws.Range("AG3").Value = 201
$ws.Range("AJ3").Value = 9
$ws.Range("AN3").Value = 8.5
$ws.Range("AY3").Value = 15

# --- Insert a new row before row 4 (the former row 4 -- Ch. Odesa vs ---
# --- Vorskla Poltava -- shifts down to row 5, unchanged) ---
$ws.Rows("4:4").Insert()

# --- Populate the newly inserted row 4 with the new fixture's data ---
# (B4 is forced to text first so Excel doesn't auto-convert the
#  "08/11/2024" literal into a date serial; the format is reset back to
#  the default "Normal" style afterwards so no stray style is left on
#  the cell.)
$ws.Range("A4").Value  = "nuwgfcy3"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value  = "08/11/2024"
$ws.Range("B4").Style  = "Normal"
$ws.Range("C4").Value  = "11:45"
$ws.Range("D4").Value  = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E4").Value  = "Al Hilal"
$ws.Range("F4").Value  = "Al Ettifaq"
$ws.Range("G4").Value  = 1.17
$ws.Range("H4").Value  = 7.5
$ws.Range("I4").Value  = 12
$ws.Range("J4").Value  = 1.5
$ws.Range("K4").Value  = 2.88
$ws.Range("L4").Value  = 9.5
$ws.Range("M4").Value  = 1.01
$ws.Range("N4").Value  = 13
$ws.Range("O4").Value  = 1.1
$ws.Range("P4").Value  = 6.5
$ws.Range("Q4").Value  = 1.36
$ws.Range("R4").Value  = 3
$ws.Range("S4").Value  = 1.22
$ws.Range("T4").Value  = 4
$ws.Range("U4").Value  = 2.1
$ws.Range("V4").Value  = 1.67
$ws.Range("W4").Value  = 10
$ws.Range("X4").Value  = 7
$ws.Range("Y4").Value  = 11
$ws.Range("Z4").Value  = 7
$ws.Range("AA4").Value = 11
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 21
$ws.Range("AD4").Value = 15
$ws.Range("AE4").Value = 29
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 900
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 51
$ws.Range("AJ4").Value = 34
$ws.Range("AK4").Value = 151
$ws.Range("AL4").Value = 81
$ws.Range("AM4").Value = 67
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 5
$ws.Range("AP4").Value = 17
$ws.Range("AQ4").Value = 11
$ws.Range("AR4").Value = 34
$ws.Range("AS4").Value = 101
$ws.Range("AT4").Value = 4
$ws.Range("AU4").Value = 11
$ws.Range("AV4").Value = 51
$ws.Range("AW4").Value = 13
$ws.Range("AX4").Value = 51
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 151
$ws.Range("BB4").Value = 500
$ws.Range("BC4").Value = 81
$ws.Range("BD4").Value = 81
